$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

$beforeCount = $tr.Paragraphs().Count

$newText = "`rDifferent data layouts entail different workload and ingestion capabilities…"
$newText += "`rBut nowadays data-intensive applications often entail both characteristics…"
$newText += "`r… Internet-of-Things networks"
$newText += "`r… Digital Twins applications"
$newText += "`r… Pervasive computing"
$newText += "`rWhat about an hybrid data structure ?"

[void]$tr.InsertAfter($newText)

$levels = @(1, 1, 2, 2, 2, 1)
$boldFlags = @($false, $false, $false, $false, $false, $true)

for ($i = 0; $i -lt $levels.Length; $i++) {
    $paraIndex = $beforeCount + 1 + $i
    $para = $tr.Paragraphs($paraIndex, 1)
    $para.IndentLevel = $levels[$i]
    if ($boldFlags[$i]) {
        $para.Font.Bold = $true
    }
}
